$d = $word.ActiveDocument

# Locate the last paragraph in the document ("Pegamos fontes open source do
# google e colocamos no nosso código.") and append two new list paragraphs
# after it, matching the existing "Trabalhando com fontes alternativas:" /
# body-text pattern used throughout the document.

$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)

# --- New paragraph 1: "Modificando o tamanho das divs:" (list level 1) ---
$endRange.InsertParagraphAfter()
$titlePara = $d.Paragraphs.Last
$titlePara.Range.ListFormat.ListLevelNumber = 2

$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertAfter(" Modificando o tamanho das divs:")

# --- New paragraph 2: body text (list level 2) ---
$titleRange = $titlePara.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()
$bodyPara = $d.Paragraphs.Last
$bodyPara.Range.ListFormat.ListLevelNumber = 3

$bodyRange = $bodyPara.Range
$bodyRange.Collapse(0)
$bodyRange.InsertAfter("Podemos alterar os tamanhos delas e o quanto delas seus itens ocupam, tudo olhando na documentação.")
